$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "GanttChart": review progress update (Part 3, 4, 5 done) + scroll the
# "Display Week" control from 10 to 5.
# ---------------------------------------------------------------------------
$gantt = $wb.Worksheets.Item("GanttChart")
$gantt.Activate()

# Move the "Display Week" scrollbar-linked cell back to week 5.
$gantt.Range("H4").Value = 5

# % DONE column: Part 3, Part 4 finished (100%), Part 5 at 75%.
$gantt.Range("H19").Value = 1
$gantt.Range("H20").Value = 1
$gantt.Range("H21").Value = 0.75

# Restore the frozen pane (rows 1-7 stay frozen) while scrolling the
# worksheet further down and updating the zoom level + selection.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$gantt.Range("A8").Select()
$win.FreezePanes = $true
$win.Zoom = 55
$win.ScrollRow = 10
$win.ScrollColumn = 1
$gantt.Range("BO16").Select()

# ---------------------------------------------------------------------------
# Sheet "Project Budget": just a scroll/selection change.
# ---------------------------------------------------------------------------
$budget = $wb.Worksheets.Item("Project Budget")
$budget.Activate()
$bwin = $excel.ActiveWindow
$bwin.ScrollRow = 24
$bwin.ScrollColumn = 2
$budget.Range("J39").Select()

# ---------------------------------------------------------------------------
# Restore GanttChart as the active/selected sheet.
# ---------------------------------------------------------------------------
$gantt.Activate()
